# The deck currently uses the "Integral" theme (ppt/theme/theme1.xml) for
# its slide master. This swaps the palette over to the stock PowerPoint
# "Office Theme" color scheme (the same 12 theme colors used by
# ppt/theme/theme2.xml, which backs the notes master).
#
# PowerPoint's ColorFormat.RGB is a plain OLE_COLOR (0x00BBGGRR), so build
# each value from its R/G/B bytes instead of hand-reversing hex strings.
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Index -> theme color slot (verified against the OOXML clrScheme order):
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
$colorScheme.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$colorScheme.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colorScheme.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$colorScheme.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colorScheme.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colorScheme.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$colorScheme.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colorScheme.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$colorScheme.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$colorScheme.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$colorScheme.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$colorScheme.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
